# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.031.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.44%  '

$ws.Range("D3").Value = "'3.512.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'594.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").Value = "'173.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.43%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("E8").Value = '  +2.56%  '

$ws.Range("E9").Value = '  +6.43%  '

$ws.Range("D10").Value = "'7.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("E11").Value = '  -1.67%  '

$ws.Range("D12").Value = "'4.118.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("D13").Value = "'0.134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("D14").Value = "'28.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.15%  '

$ws.Range("E15").Value = '  -0.02%  '

$ws.Range("D16").Value = "'67.028.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").Value = "'3.515.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").Value = "'6.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.47%  '

$ws.Range("D19").Value = "'14.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("D20").Value = "'394.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.45%  '

$ws.Range("D21").Value = "'7.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("D22").Value = "'73.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.17%  '

$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("E24").Value = '  +1.30%  '

$ws.Range("E25").Value = '  -4.94%  '

$ws.Range("D26").Value = "'10.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.29%  '

$ws.Range("D27").Value = "'0.180"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.29%  '

$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.17%  '

$ws.Range("E29").Value = '  -2.55%  '

$ws.Range("E30").Value = '  -2.20%  '

$ws.Range("E31").Value = '  -0.11%  '

$ws.Range("D32").Value = "'23.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.18%  '

$ws.Range("D33").Value = "'7.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.31%  '

$ws.Range("E34").Value = '  +2.76%  '

$ws.Range("D35").Value = "'163.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.86%  '

$ws.Range("E36").Value = '  -0.85%  '

$ws.Range("E37").Value = '  -2.02%  '

$ws.Range("D38").Value = "'6.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.98%  '

$ws.Range("D39").Value = "'4.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("D40").Value = "'0.0743"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.45%  '

$ws.Range("D41").Value = "'27.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("E42").Value = '  -1.21%  '

$ws.Range("D43").Value = "'2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("D44").Value = "'2.800.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("D45").Value = "'42.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.45%  '

$ws.Range("D46").Value = "'0.0305"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.41%  '

$ws.Range("D47").Value = "'341.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.92%  '

$ws.Range("D48").Value = "'1.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.27%  '

$ws.Range("D49").Value = "'33.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.08%  '

$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("D51").Value = "'0.847"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.61%  '
